# Scheduled-runner price/profit refresh for the Chocobo_Profits workbook.
# Updates the cached market-board figures (currentAveragePrice*, Leve
# price/profit columns H-N) for a handful of rows across the eight job
# sheets, matching the latest API pull.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 10004976
$ws.Range("I74").Value = 33336652
$ws.Range("J74").Value = 5685.7144
$ws.Range("K74").Value = 33336652
$ws.Range("L74").Value = 5685.7144
$ws.Range("M74").Value = -33335716
$ws.Range("N74").Value = -7557.7144

$ws.Range("H76").Value = 3207.739
$ws.Range("I76").Value = 3013.7
$ws.Range("K76").Value = 3013.7
$ws.Range("M76").Value = -2698.7

$ws.Range("H77").Value = 10004976
$ws.Range("I77").Value = 33336652
$ws.Range("J77").Value = 5685.7144
$ws.Range("K77").Value = 166683260
$ws.Range("L77").Value = 28428.572
$ws.Range("M77").Value = -166678580
$ws.Range("N77").Value = -37788.572

$ws.Range("H79").Value = 3207.739
$ws.Range("I79").Value = 3013.7
$ws.Range("K79").Value = 3013.7
$ws.Range("M79").Value = -1921.7

$ws.Range("H98").Value = 2286.8293
$ws.Range("I98").Value = 728.05884
$ws.Range("J98").Value = 9858
$ws.Range("K98").Value = 728.05884
$ws.Range("L98").Value = 9858
$ws.Range("M98").Value = 769.94116
$ws.Range("N98").Value = -12854

$ws.Range("H113").Value = 6843.7144
$ws.Range("I113").Value = 2950
$ws.Range("J113").Value = 8401.200000000001
$ws.Range("K113").Value = 2950
$ws.Range("L113").Value = 8401.200000000001
$ws.Range("M113").Value = 304
$ws.Range("N113").Value = -14909.2

$ws.Range("H122").Value = 2286.8293
$ws.Range("I122").Value = 728.05884
$ws.Range("J122").Value = 9858
$ws.Range("K122").Value = 2184.17652
$ws.Range("L122").Value = 29574
$ws.Range("M122").Value = 265.82348
$ws.Range("N122").Value = -34474

$ws.Range("H137").Value = 4415.727
$ws.Range("I137").Value = 2259.3635
$ws.Range("K137").Value = 6778.0905
$ws.Range("M137").Value = -4228.0905

# ---------------------------------------------------------------- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5428.698
$ws.Range("I32").Value = 4302
$ws.Range("J32").Value = 9034.134
$ws.Range("K32").Value = 4302
$ws.Range("L32").Value = 9034.134
$ws.Range("M32").Value = -4015
$ws.Range("N32").Value = -9608.134

$ws.Range("H74").Value = 4585.4
$ws.Range("I74").Value = 4795.7896
$ws.Range("J74").Value = 3919.1667
$ws.Range("K74").Value = 4795.7896
$ws.Range("L74").Value = 3919.1667
$ws.Range("M74").Value = -3921.7896
$ws.Range("N74").Value = -5667.1667

$ws.Range("H77").Value = 4585.4
$ws.Range("I77").Value = 4795.7896
$ws.Range("J77").Value = 3919.1667
$ws.Range("K77").Value = 23978.948
$ws.Range("L77").Value = 19595.8335
$ws.Range("M77").Value = -19610.948
$ws.Range("N77").Value = -28331.8335

$ws.Range("H88").Value = 6670331.5
$ws.Range("I88").Value = 8336664
$ws.Range("K88").Value = 8336664
$ws.Range("M88").Value = -8336258

$ws.Range("H91").Value = 6670331.5
$ws.Range("I91").Value = 8336664
$ws.Range("K91").Value = 8336664
$ws.Range("M91").Value = -8335260

$ws.Range("H97").Value = 993.9474
$ws.Range("I97").Value = 982.5
$ws.Range("K97").Value = 982.5
$ws.Range("M97").Value = -486.5

$ws.Range("H132").Value = 3075.6086
$ws.Range("I132").Value = 1617.2858
$ws.Range("J132").Value = 5344.1113
$ws.Range("K132").Value = 4851.857400000001
$ws.Range("L132").Value = 16032.3339
$ws.Range("M132").Value = -2321.857400000001
$ws.Range("N132").Value = -21092.3339

$ws.Range("H139").Value = 43515
$ws.Range("J139").Value = 43515
$ws.Range("L139").Value = 43515
$ws.Range("N139").Value = -53795

# ---------------------------------------------------------------- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H74").Value = 45400
$ws.Range("J74").Value = 45400
$ws.Range("L74").Value = 45400
$ws.Range("N74").Value = -47272

$ws.Range("H77").Value = 45400
$ws.Range("J77").Value = 45400
$ws.Range("L77").Value = 136200
$ws.Range("N77").Value = -145560

$ws.Range("H86").Value = 2146.6924
$ws.Range("I86").Value = 1887.5
$ws.Range("J86").Value = 2561.4
$ws.Range("K86").Value = 1887.5
$ws.Range("L86").Value = 2561.4
$ws.Range("M86").Value = -764.5
$ws.Range("N86").Value = -4807.4

$ws.Range("H89").Value = 2146.6924
$ws.Range("I89").Value = 1887.5
$ws.Range("J89").Value = 2561.4
$ws.Range("K89").Value = 9437.5
$ws.Range("L89").Value = 12807
$ws.Range("M89").Value = -3821.5
$ws.Range("N89").Value = -24039

$ws.Range("H94").Value = 1339.875
$ws.Range("I94").Value = 1534.8334
$ws.Range("J94").Value = 755
$ws.Range("K94").Value = 1534.8334
$ws.Range("L94").Value = 755
$ws.Range("M94").Value = -1083.8334
$ws.Range("N94").Value = -1657

$ws.Range("H99").Value = 1988.0625
$ws.Range("I99").Value = 1268.7
$ws.Range("J99").Value = 3187
$ws.Range("K99").Value = 1268.7
$ws.Range("L99").Value = 3187
$ws.Range("M99").Value = 229.3
$ws.Range("N99").Value = -6183

# ---------------------------------------------------------------- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11907158
$ws.Range("I31").Value = 1336.9259
$ws.Range("J31").Value = 33337636
$ws.Range("K31").Value = 1336.9259
$ws.Range("L31").Value = 33337636
$ws.Range("M31").Value = -1041.9259
$ws.Range("N31").Value = -33338226

$ws.Range("H34").Value = 11907158
$ws.Range("I34").Value = 1336.9259
$ws.Range("J34").Value = 33337636
$ws.Range("K34").Value = 1336.9259
$ws.Range("L34").Value = 33337636
$ws.Range("M34").Value = -1134.9259
$ws.Range("N34").Value = -33338040

$ws.Range("H132").Value = 1233.7119
$ws.Range("I132").Value = 704.4773
$ws.Range("J132").Value = 2786.1333
$ws.Range("K132").Value = 2113.4319
$ws.Range("L132").Value = 8358.3999
$ws.Range("M132").Value = 416.5681
$ws.Range("N132").Value = -13418.3999

# ---------------------------------------------------------------- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 43000000
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
# Row 33 no longer has an HQ leve (N33 cell is removed outright, not zeroed).
$ws.Range("N33").ClearContents()

$ws.Range("H70").Value = 6288.0312
$ws.Range("I70").Value = 5789.115
$ws.Range("J70").Value = 8450
$ws.Range("K70").Value = 5789.115
$ws.Range("L70").Value = 8450
$ws.Range("M70").Value = -5519.115
$ws.Range("N70").Value = -8990

$ws.Range("H73").Value = 6288.0312
$ws.Range("I73").Value = 5789.115
$ws.Range("J73").Value = 8450
$ws.Range("K73").Value = 5789.115
$ws.Range("L73").Value = 8450
$ws.Range("M73").Value = -4853.115
$ws.Range("N73").Value = -10322

$ws.Range("H139").Value = 96666.664
$ws.Range("J139").Value = 96666.664
$ws.Range("L139").Value = 96666.664
$ws.Range("N139").Value = -106946.664

$ws.Range("H140").Value = 38702
$ws.Range("J140").Value = 38702
$ws.Range("L140").Value = 38702
$ws.Range("N140").Value = -49062

# ---------------------------------------------------------------- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H111").Value = 39750
$ws.Range("J111").Value = 39750
$ws.Range("L111").Value = 39750
$ws.Range("N111").Value = -47930

$ws.Range("H132").Value = 10092.071
$ws.Range("I132").Value = 10930.96
$ws.Range("J132").Value = 8858.412
$ws.Range("K132").Value = 32792.88
$ws.Range("L132").Value = 26575.236
$ws.Range("M132").Value = -30262.88
$ws.Range("N132").Value = -31635.236

$ws.Range("H139").Value = 40254
$ws.Range("J139").Value = 40254
$ws.Range("L139").Value = 40254
$ws.Range("N139").Value = -50534

$ws.Range("H140").Value = 80696.09
$ws.Range("J140").Value = 80696.09
$ws.Range("L140").Value = 80696.09
$ws.Range("N140").Value = -91056.09

$ws.Range("H141").Value = 41553.41
$ws.Range("J141").Value = 41553.41
$ws.Range("L141").Value = 41553.41
$ws.Range("N141").Value = -51913.41

# ---------------------------------------------------------------- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H138").Value = 44560.89
$ws.Range("J138").Value = 44560.89
$ws.Range("L138").Value = 44560.89
$ws.Range("N138").Value = -54840.89
